# "unify the conception of DataNode, DataTable, Entity."
# Rename the single worksheet from the old "Property1" label to the
# unified "DataNode" name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Leave the selection where the authoring session left it.
$ws.Range("C41").Select() | Out-Null
